# Update the "counter" values (column D) for the RFID users table.
# alice: 17 -> 16, bob: 16 -> 14, charlie: 6 -> 5, Jacobe stays at 2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 16
$ws.Range("D3").Value = 14
$ws.Range("D4").Value = 5
